# Final fixes, balancing, and mastering.
# Apply content corrections to the AssetList sheet: fix typos ("wonce"->"once",
# "finished"->"made"), rework several Notes/Description/Assets-Required
# strings (pinging/beeping -> clicking/heartbeat/jingle wording), swap the
# placeholder "Assets Required" text for Title/Dungeon/Boss/Defeat/Victory
# Music rows with proper descriptive asset names, and reclassify the
# Game Defeat / Game Victory rows from "Sound Effect" to "Interface".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Title Music: Assets Required description reworded
$ws.Range("D2").Value = "Slow, adventurous music track"

# Row 3 - Select: pinging -> clicking wording, and "finished" -> "made" in Notes
$ws.Range("C3").Value = "A medium-pitch clicking noise that reinforces the change of a selected menu item."
$ws.Range("D3").Value = "Medium-pitch clicking noise"
$ws.Range("F3").Value = "Sound was previously made as an electronic pinging noise, but has been remade as a wooden clicking noise. "

# Row 4 - Confirm: pinging -> clicking wording, and "finished" -> "made" in Notes
$ws.Range("C4").Value = "A clikcing noise with increasing pitch that reinforces the advancement of the game."
$ws.Range("D4").Value = "Medium-pitch clicking noise"
$ws.Range("F4").Value = "Sound was previously made as an electronic pinging noise, but has been remade as a wooden clicking noise. "

# Row 5 - Dungeon Music: Assets Required description reworded
$ws.Range("D5").Value = "Underground-themed music track"

# Row 6 - Boss Music: fix "wonce" typo -> "once"; Assets Required reworded
$ws.Range("C6").Value = "A medium-paced music track that plays during the boss fight. The music changes once the boss has reached half health."
$ws.Range("D6").Value = "Dramatic orchestral music track"

# Row 7 - Game Defeat: recategorized Sound Effect -> Interface; chime -> jingle wording
$ws.Range("B7").Value = "Interface"
$ws.Range("C7").Value = "A defeated jingle that plays on the menu after the player has been defeated."
$ws.Range("D7").Value = "Defeated-sounding jingle"

# Row 8 - Game Victory: recategorized Sound Effect -> Interface; chime -> jingle wording
$ws.Range("B8").Value = "Interface"
$ws.Range("C8").Value = "A triumphant jingle that plays on the menu after the player has beaten the game."
$ws.Range("D8").Value = "Triumphant-sounding jingle"

# Row 9 - Defeat Music: "corresponding chime" -> "defeat chime"; Assets Required reworded
$ws.Range("C9").Value = "A slow-paced music track that plays on the menu once the game has been lost and the defeat chime has finished. "
$ws.Range("D9").Value = "Slow, ambient music track"

# Row 10 - Victory Music: "corresponding chime" -> "victory chime"; Assets Required reworded
$ws.Range("C10").Value = "A slow-paced music track that plays on the menu once the game has been won and the victory chime has finished. "
$ws.Range("D10").Value = "Slow, ambient music track"

# Row 11 - Player Low Health: beeping -> heartbeat wording
$ws.Range("C11").Value = " A low-pitch heartbeat that plays when the player has one bubble or less of health remaining"
$ws.Range("D11").Value = "Heartbeat sound"

# Row 15 - Player Block Projectile: "finished" -> "made" in Notes
$ws.Range("F15").Value = "Sound was previously made as a shield clinking noise, but has been remade as a magic noise. `n"

# Row 2's description got shorter (now fits on one line), so Excel's autofit
# drops the explicit row height that was sized for two wrapped lines.
$ws.Rows.Item(2).AutoFit()

# The author's cursor ended up on C7 when the file was last saved.
[void]$ws.Range("C7").Select()
